$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.196.05'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.827.09'
$ws.Range("E3").Value = '  +2.40%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''225.01'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '''0.560'
$ws.Range("E6").Value = '  +1.36%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''31.97'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +3.77%  '
$ws.Range("D10").Value = '''0.0725'
$ws.Range("E10").Value = '  +10.41%  '
$ws.Range("D11").Value = '''0.0931'
$ws.Range("D12").Value = '2.091.25'
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("D13").Value = '1.825.89'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("D14").Value = '''10.83'
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("E15").Value = '  +2.90%  '
$ws.Range("D16").Value = '34.213.14'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '''4.34'
$ws.Range("E17").Value = '  +3.17%  '
$ws.Range("D18").Value = '''69.73'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '''251.12'
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("D20").Value = '0.0₃0792'
$ws.Range("E20").Value = '  +7.04%  '
$ws.Range("D21").Value = '''11.16'
$ws.Range("E21").Value = '  +7.77%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''4.27'
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").Value = '''2.16'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = '''160.56'
$ws.Range("E25").Value = '  +2.03%  '
$ws.Range("D26").Value = '''16.67'
$ws.Range("E26").Value = '  +1.68%  '
$ws.Range("D27").Value = '''7.27'
$ws.Range("E27").Value = '  +3.66%  '
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''0.0534'
$ws.Range("E30").Value = '  +3.94%  '
$ws.Range("D31").Value = '''3.78'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '''1.21'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("E34").Value = '  +1.53%  '
$ws.Range("D35").Value = '1.438.23'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = '''0.646'
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("D39").Value = '''0.964'
$ws.Range("E39").Value = '  +8.36%  '
$ws.Range("D40").Value = '''81.78'
$ws.Range("E40").Value = '  -1.43%  '
$ws.Range("E41").Value = '  -2.97%  '
$ws.Range("D42").Value = '''2.35'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +5.19%  '
$ws.Range("D44").Value = '''6.10'
$ws.Range("E44").Value = '  +4.43%  '
$ws.Range("D45").Value = '1.986.90'
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").Value = '''107.47'
$ws.Range("E48").Value = '  +8.84%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''11.90'
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("E51").Value = '  +5.00%  '
